# Replace the placeholder tokens in column A (rows 2-10) with their real
# values, matching the author's "not finished but last day" commit.
#
# A2 ("202004142359.") looks numeric, so Excel would otherwise silently
# coerce it to a number and drop the trailing period; force the cell to
# Text first, assign, then clear the format override back off so the
# cell keeps the workbook's original (unstyled) look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "202004142359."
$ws.Range("A3").Value = "TOP."
$ws.Range("A4").Value = "prj-isp-t4777-appl-svil-001.T4777W.VD_T4_DB_PSET.NOM_PSET"
$ws.Range("A5").Value = "prj-isp-t4777-appl-svil-001.T4777W.VD_T4_DB_PSET.COD_PSET"
$ws.Range("A6").Value = "prj-isp-t4777-appl-svil-001.T4777W.VD_T4_DB_PSET.NOM_PSET"
$ws.Range("A7").Value = "NUM_AA_TT"
$ws.Range("A8").Value = " COD_KEY_PIAZZA_REGOLAMENTO"
$ws.Range("A9").Value = "TOP"
$ws.Range("A10").Value = "T4"
$ws.Range("A2").ClearFormats()

# Narrow the two data columns a bit (A: ~66.36 -> ~55.73, B: ~33.09 -> ~30.18
# characters).
$ws.Columns.Item(1).ColumnWidth = 54.83333333333333
$ws.Columns.Item(2).ColumnWidth = 29.333333333333336
